# Slide 7 contains the "Table with Cell Styling" demo table (shape 2).
# Every run in the table gains an explicit i="0" (Italic = False) attribute,
# and the runs that previously had no b attribute at all gain an explicit
# b="0" (Bold = False) attribute as well. Runs that were already explicitly
# bold (b="1") keep that attribute untouched.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        $tr = $cell.Shape.TextFrame.TextRange

        # Only ever assign False here - re-asserting an already-true Bold
        # value corrupts the Bold reads of other runs in this runtime, so
        # cells that are already bold are simply left alone.
        $wasBold = $tr.Font.Bold
        if ($wasBold -eq 0) {
            $tr.Font.Bold = 0
        }

        $tr.Font.Italic = 0
    }
}
